$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price cells whose new values look like plain numbers,
# so Excel stores them as literal text strings (matching the source data format)
# instead of auto-converting them into numeric values.
$numericLookingCells = @("D4", "D5", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Updated Price (D) and Volume(1h) (E) values per the latest cryptos snapshot
$ws.Range("D2").Value = '25.881.79'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.729.87'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '244.70'
$ws.Range("E5").Value = '  +2.75%  '
$ws.Range("D7").Value = '0.5019'
$ws.Range("E7").Value = '  -3.07%  '
$ws.Range("D8").Value = '0.2706'
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("D9").Value = '0.06150'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '1.731.59'
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").Value = '0.07223'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '15.10'
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("D13").Value = '0.6467'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '4.733'
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("D15").Value = '76.72'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '25.884.63'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '11.87'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '0.000006800'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '4.584'
$ws.Range("E21").Value = '  +7.05%  '
$ws.Range("D22").Value = '1.956.91'
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").Value = '8.769'
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").Value = '5.460'
$ws.Range("E24").Value = '  +3.94%  '
$ws.Range("D25").Value = '134.26'
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("D26").Value = '15.26'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").Value = '1.414'
$ws.Range("E27").Value = '  -6.58%  '
$ws.Range("D28").Value = '1.778'
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("D29").Value = '105.06'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").Value = '3.951'
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").Value = '0.08109'
$ws.Range("E31").Value = '  -2.24%  '
$ws.Range("D32").Value = '3.686'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").Value = '0.04707'
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("D34").Value = '2.652'
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").Value = '0.9962'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").Value = '0.6093'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("D37").Value = '2.740'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("D38").Value = '0.01605'
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '0.8676'
$ws.Range("E39").Value = '  +17.24%  '
$ws.Range("D40").Value = '1.946'
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").Value = '1.000'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '101.51'
$ws.Range("E42").Value = '  +3.71%  '
$ws.Range("D43").Value = '0.3872'
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = '4.992'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '0.1180'
$ws.Range("E45").Value = '  +4.43%  '
$ws.Range("D46").Value = '6.329'
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").Value = '55.52'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = '30.68'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value = '7.656'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("D51").Value = '0.3464'
$ws.Range("E51").Value = '  +1.71%  '
